# Update to latest versions.
#
# This script:
#  1) Fixes the "Find_SNP_groups.py" label (capitalization typo) on slide 7
#     to "find_SNP_groups.py".
#  2) Duplicates the whole 13-shape flow-chart diagram on slide 7 and moves
#     the copy down/right, to show the second ("kanika") pipeline branch.
#  3) Re-labels two of the duplicated shapes for the new branch:
#       "snv_plus_indels.20180919.csv" -> "lucian_from_kanika.csv"
#       "VAFclusters/"                 -> "VAFclusters_kanika/"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- 1) fix the capitalization of the existing "Find_SNP_groups.py" oval ---
$findOval = $s.Shapes.Item(11)   # Oval 61, id=62
$findOval.TextFrame.TextRange.Text = "find_SNP_groups.py"

# --- 2) duplicate every shape in the diagram and move the copies ---
# Uniform translation applied to the whole diagram (in points; 12700 EMU = 1pt)
$dx = 152400 / 12700      # = 12 pt
$dy = 3634494 / 12700     # ~= 286.18 pt

# Target absolute positions/sizes (points) for each duplicate, taken from the
# final layout of the second copy of the diagram.
# NOTE: left/top/width/height below are deliberately given at float32
# precision (rather than the nearest decimal) because this host stores
# shape geometry as 32-bit floats internally; using the closest decimal
# value can round to the EMU neighbor when re-quantized. These values are
# the float32 representable numbers that round-trip to the exact target
# EMU (1 pt = 12700 EMU) offsets/extents of the duplicated diagram.
$targets = @(
    @{ idx = 1;  left = 438.8251647949219; top = 480.8587951660156; width = 187.5821685791016; height = 26.6577568054199 },  # TextBox 1
    @{ idx = 2;  left = 221.3016204833984; top = 337.0431213378906; width = 273.9235229492188;  height = 36.918701171875 },   # Oval 3
    @{ idx = 3;  left = 345.36328125;      top = 380.8620300292969; width = 19.8003559112549;   height = 5.9999904632568 },  # Elbow Connector 47
    @{ idx = 4;  left = 371.9999694824219; top = 427.4998168945312; width = 273.9235229492188;  height = 36.918701171875 },  # Oval 12
    @{ idx = 5;  left = 512.5689697265625; top = 460.8114013671875; width = 16.4401988983154;   height = 23.6544494628906 }, # Elbow Connector 18 (a)
    @{ idx = 6;  left = 199.5925598144531; top = 291.5216369628906; width = 317.341552734375;   height = 29.0812225341797 }, # TextBox 15
    @{ idx = 7;  left = 350.0432739257812; top = 328.8228759765625; width = 16.4401988983154;   height = 0.0000787400131 },  # Elbow Connector 30
    @{ idx = 8;  left = 277.4201354980469; top = 393.7622680664062; width = 149.6865081787109;  height = 26.6577568054199 }, # TextBox 19
    @{ idx = 9;  left = 427.0727233886719; top = 345.6106872558594; width = 7.0797543525696;    height = 156.6983184814453 }, # Elbow Connector 41 (a)
    @{ idx = 10; left = 72.1694946289062;  top = 479.2373046875;    width = 143.1524810791016;  height = 26.6577568054199 }, # TextBox 60
    @{ idx = 11; left = 31.4619312286377;  top = 427.4998168945312; width = 273.9235229492188;  height = 36.918701171875 },  # Oval 61
    @{ idx = 12; left = 148.6753997802734; top = 459.4889526367188; width = 14.8187208175659;   height = 24.6779136657715 }, # Elbow Connector 18 (b)
    @{ idx = 13; left = 256.8037414550781; top = 332.0400695800781; width = 7.0797543525696;    height = 183.8396606445312 } # Elbow Connector 41 (b)
)

$origCount = $s.Shapes.Count
$dupMap = @{}

for ($i = 1; $i -le $origCount; $i++) {
    $orig = $s.Shapes.Item($i)
    $dup = $orig.Duplicate()
    $dupMap[$i] = $dup
}

foreach ($t in $targets) {
    $dup = $dupMap[$t.idx]
    $dup.Left = $t.left
    $dup.Top = $t.top
    $dup.Width = $t.width
    $dup.Height = $t.height
}

# --- 3) re-label the two shapes that differ in the new ("kanika") branch ---
$dupCsvLabel = $dupMap[6]     # duplicate of TextBox 15 ("snv_plus_indels.20180919.csv")
$dupCsvLabel.TextFrame.TextRange.Text = "lucian_from_kanika.csv"

$dupFolderLabel = $dupMap[8]  # duplicate of TextBox 19 ("VAFclusters/")
$dupFolderLabel.TextFrame.TextRange.Text = "VAFclusters_kanika/"
